# Applies the edits described in commit "corregir un par de detalles"
# to the "Una sugerencia practica de versionado semantico" document.

$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) "=> Aumentamos X cuando se hace un lanzamiento completo a
#    "producción"." -> replace "hace" with "emite"
# ----------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "Aumentamos X cuando se hace un lanzamiento completo a", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "Aumentamos X cuando se emite un lanzamiento completo a", 2)
if (-not $found) { Write-Host "WARN: replacement 1 not found" }

# ----------------------------------------------------------------------
# 2) "=> Aumentamos Y cuando se lanza una corrección a la versión X."
#    -> "=> Aumentamos Y cuando se emite una corrección sobre la
#    última versión X. en producción"
# ----------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "lanza una corrección a la versión X.", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "emite una corrección sobre la última versión X. en producción", 2)
if (-not $found) { Write-Host "WARN: replacement 2 not found" }

# ----------------------------------------------------------------------
# 3) Drop the trailing parenthetical note after "X.Y.3.0, etc."
# ----------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    " (nota: se pueden abreviar a X.Y.1, X.Y.2, X.Y.3, etc.)", $true,
    $false, $false, $false, $false, $true, 1, $false, "", 2)
if (-not $found) { Write-Host "WARN: replacement 3 not found" }

# ----------------------------------------------------------------------
# 4) "emite una versión para pruebas" -> "emite una versión “para
#    pruebas” internas"
# ----------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "emite una versión para pruebas", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "emite una versión " + [char]0x201C + "para pruebas" + [char]0x201D + " internas", 2)
if (-not $found) { Write-Host "WARN: replacement 4 not found" }

# ----------------------------------------------------------------------
# 5) "empezamos a trabajar después de esa emisión." -> "empezamos a
#    trabajar después de cualquier emisión."
# ----------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "empezamos a trabajar después de esa emisión.", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "empezamos a trabajar después de cualquier emisión.", 2)
if (-not $found) { Write-Host "WARN: replacement 5 not found" }

# ----------------------------------------------------------------------
# 6) Table column widths shift by 1 twip (4704->4703, 2377->2378) as a
#    result of Word's layout recalculation. Nudge the first and last
#    columns accordingly (columns are expressed in points = twips/20).
# ----------------------------------------------------------------------
$tbl = $d.Tables.Item(1)
$tbl.Columns.Item(1).Width = 235.15
$tbl.Columns.Item(4).Width = 118.9

Write-Host "Done"
